# Weekly update: insert two new price rows (new week's data) right before
# the current row 205, shifting all the existing rows 205-340 down by two
# rows (they become rows 207-342).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 205 (inserting twice at the same
# address pushes everything down by 2).
$ws.Range("A205").EntireRow.Insert()
$ws.Range("A205").EntireRow.Insert()

# --- New row 205 ---
$ws.Cells.Item(205, 1).Value  = 10
$ws.Cells.Item(205, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(205, 3).Value  = "La Araucanía"
$ws.Cells.Item(205, 4).Value  = 44767
$ws.Cells.Item(205, 5).Value  = 9
$ws.Cells.Item(205, 6).Value  = 100112044
$ws.Cells.Item(205, 7).Value  = "Perejil"
$ws.Cells.Item(205, 8).Value  = "Sin especificar"
$ws.Cells.Item(205, 9).Value  = "Primera"
$ws.Cells.Item(205, 10).Value = 50
$ws.Cells.Item(205, 11).Value = 6000
$ws.Cells.Item(205, 12).Value = 6000
$ws.Cells.Item(205, 13).Value = 6000
$ws.Cells.Item(205, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(205, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(205, 16).Value = 2000
$ws.Cells.Item(205, 17).Value = 3
$ws.Cells.Item(205, 18).Value = "Hortaliza"

# --- New row 206 ---
$ws.Cells.Item(206, 1).Value  = 10
$ws.Cells.Item(206, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(206, 3).Value  = "La Araucanía"
$ws.Cells.Item(206, 4).Value  = 44767
$ws.Cells.Item(206, 5).Value  = 9
$ws.Cells.Item(206, 6).Value  = 100112044
$ws.Cells.Item(206, 7).Value  = "Perejil"
$ws.Cells.Item(206, 8).Value  = "Sin especificar"
$ws.Cells.Item(206, 9).Value  = "Primera"
$ws.Cells.Item(206, 10).Value = 40
$ws.Cells.Item(206, 11).Value = 4300
$ws.Cells.Item(206, 12).Value = 4300
$ws.Cells.Item(206, 13).Value = 4300
$ws.Cells.Item(206, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(206, 15).Value = "Región Metropolitana"
$ws.Cells.Item(206, 16).Value = 1433
$ws.Cells.Item(206, 17).Value = 3
$ws.Cells.Item(206, 18).Value = "Hortaliza"

# Make sure the date cells keep the same date number format as the rest
# of column D.
$ws.Range("D205").NumberFormat = $ws.Range("D204").NumberFormat
$ws.Range("D206").NumberFormat = $ws.Range("D204").NumberFormat
